$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.384.89'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.716.59'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9969'
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.47'
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4867'
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2588'
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06170'
$ws.Range("E9").Value = '  -2.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.721.67'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06950'
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.50'
$ws.Range("E12").Value = '  -1.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.491'
$ws.Range("E13").Value = '  -2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5979'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.59'
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9981'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.370.37'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9969'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("E19").Value = '  -4.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.25'
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.946.50'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.418'
$ws.Range("E22").Value = '  -3.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.454'
$ws.Range("E23").Value = '  -2.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.063'
$ws.Range("E24").Value = '  -3.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.84'
$ws.Range("E25").Value = '  -2.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.24'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.407'
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.737'
$ws.Range("E28").Value = '  -1.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.87'
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07960'
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.612'
$ws.Range("E32").Value = '  -2.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04442'
$ws.Range("E33").Value = '  -2.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.600'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9957'
$ws.Range("E35").Value = '  -1.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6185'
$ws.Range("E36").Value = '  -2.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9364'
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.976'
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.372'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9973'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01475'
$ws.Range("E41").Value = '  -2.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.12'
$ws.Range("E42").Value = '  -2.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.448'
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3812'
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.831'
$ws.Range("E45").Value = '  -0.89%  '
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05349'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.39'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.716'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.23'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.217'
$ws.Range("E51").Value = '  -3.50%  '
